$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (1215) down to the new rows (1216-1229)
$fmtSrc = $ws.Range("A1215:V1215")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A1216:V1229")
$fmtDst.PasteSpecial(-4122)

# Row 1216 column E (new player "Nathanael Beta") uses the highlighted name style (s=6),
# matching the style already used elsewhere in the sheet for newly-added players (e.g. E993).
$eStyleSrc = $ws.Range("E993")
$eStyleSrc.Copy()
$eStyleDst = $ws.Range("E1216")
$eStyleDst.PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 1216
$ws.Range("A1216").Value = "Entrainement"
$ws.Range("B1216").Value = 46050
$ws.Range("C1216").Value = "Global"
$ws.Range("D1216").Value = "J-3"
$ws.Range("E1216").Value = "Nathanael Beta"
$ws.Range("F1216").Value = "left back"
$ws.Range("G1216").Value = "01:29:49"
$ws.Range("H1216").Value = 6.46
$ws.Range("I1216").Value = 0.89
$ws.Range("J1216").Value = 5.56
$ws.Range("K1216").Value = 0.51
$ws.Range("L1216").Value = 0.34
$ws.Range("M1216").Value = 0.05
$ws.Range("N1216").Value = 0
$ws.Range("O1216").Value = 5
$ws.Range("P1216").Value = 4.25
$ws.Range("Q1216").Value = 28.45
$ws.Range("R1216").Value = 4.11
$ws.Range("S1216").Value = 42
$ws.Range("T1216").Value = 3
$ws.Range("U1216").Value = 13
$ws.Range("V1216").Value = 5

# Row 1217
$ws.Range("A1217").Value = "Entrainement"
$ws.Range("B1217").Value = 46050
$ws.Range("C1217").Value = "Global"
$ws.Range("D1217").Value = "J-3"
$ws.Range("E1217").Value = "Yoan Zouma"
$ws.Range("F1217").Value = "center back"
$ws.Range("G1217").Value = "01:39:55"
$ws.Range("H1217").Value = 6.36
$ws.Range("I1217").Value = 0.4
$ws.Range("J1217").Value = 5.95
$ws.Range("K1217").Value = 0.37
$ws.Range("L1217").Value = 0.03
$ws.Range("M1217").Value = 0.01
$ws.Range("N1217").Value = 0
$ws.Range("O1217").Value = 1
$ws.Range("P1217").Value = 3.79
$ws.Range("Q1217").Value = 27.46
$ws.Range("R1217").Value = 3.92
$ws.Range("S1217").Value = 19
$ws.Range("T1217").Value = 0
$ws.Range("U1217").Value = 6
$ws.Range("V1217").Value = 1

# Row 1218
$ws.Range("A1218").Value = "Entrainement"
$ws.Range("B1218").Value = 46050
$ws.Range("C1218").Value = "Global"
$ws.Range("D1218").Value = "J-3"
$ws.Range("E1218").Value = "Karahali Souaré"
$ws.Range("F1218").Value = "right forward"
$ws.Range("G1218").Value = "01:39:11"
$ws.Range("H1218").Value = 8.4
$ws.Range("I1218").Value = 1.37
$ws.Range("J1218").Value = 6.99
$ws.Range("K1218").Value = 0.77
$ws.Range("L1218").Value = 0.49
$ws.Range("M1218").Value = 0.14
$ws.Range("N1218").Value = 0
$ws.Range("O1218").Value = 15
$ws.Range("P1218").Value = 4.84
$ws.Range("Q1218").Value = 30.04
$ws.Range("R1218").Value = 5.25
$ws.Range("S1218").Value = 97
$ws.Range("T1218").Value = 37
$ws.Range("U1218").Value = 78
$ws.Range("V1218").Value = 28

# Row 1219
$ws.Range("A1219").Value = "Entrainement"
$ws.Range("B1219").Value = 46050
$ws.Range("C1219").Value = "Global"
$ws.Range("D1219").Value = "J-3"
$ws.Range("E1219").Value = "Ilan Ihaddadene"
$ws.Range("F1219").Value = "center midfield"
$ws.Range("G1219").Value = "01:40:56"
$ws.Range("H1219").Value = 7.84
$ws.Range("I1219").Value = 0.8
$ws.Range("J1219").Value = 7.03
$ws.Range("K1219").Value = 0.59
$ws.Range("L1219").Value = 0.2
$ws.Range("M1219").Value = 0.03
$ws.Range("N1219").Value = 0
$ws.Range("O1219").Value = 3
$ws.Range("P1219").Value = 4.61
$ws.Range("Q1219").Value = 27.77
$ws.Range("R1219").Value = 4.54
$ws.Range("S1219").Value = 35
$ws.Range("T1219").Value = 3
$ws.Range("U1219").Value = 14
$ws.Range("V1219").Value = 4

# Row 1220
$ws.Range("A1220").Value = "Entrainement"
$ws.Range("B1220").Value = 46050
$ws.Range("C1220").Value = "Global"
$ws.Range("D1220").Value = "J-3"
$ws.Range("E1220").Value = "Naim Ighbane"
$ws.Range("F1220").Value = "center back"
$ws.Range("G1220").Value = "01:38:51"
$ws.Range("H1220").Value = 7.23
$ws.Range("I1220").Value = 0.51
$ws.Range("J1220").Value = 6.71
$ws.Range("K1220").Value = 0.37
$ws.Range("L1220").Value = 0.11
$ws.Range("M1220").Value = 0.04
$ws.Range("N1220").Value = 0
$ws.Range("O1220").Value = 5
$ws.Range("P1220").Value = 3.89
$ws.Range("Q1220").Value = 27.94
$ws.Range("R1220").Value = 4.44
$ws.Range("S1220").Value = 15
$ws.Range("T1220").Value = 3
$ws.Range("U1220").Value = 26
$ws.Range("V1220").Value = 1

# Row 1221
$ws.Range("A1221").Value = "Entrainement"
$ws.Range("B1221").Value = 46050
$ws.Range("C1221").Value = "Global"
$ws.Range("D1221").Value = "J-3"
$ws.Range("E1221").Value = "Sofiane Belle"
$ws.Range("F1221").Value = "left forward"
$ws.Range("G1221").Value = "01:39:09"
$ws.Range("H1221").Value = 5.42
$ws.Range("I1221").Value = 0.38
$ws.Range("J1221").Value = 5.03
$ws.Range("K1221").Value = 0.24
$ws.Range("L1221").Value = 0.12
$ws.Range("M1221").Value = 0.02
$ws.Range("N1221").Value = 0
$ws.Range("O1221").Value = 2
$ws.Range("P1221").Value = 3.16
$ws.Range("Q1221").Value = 29.19
$ws.Range("R1221").Value = 4.36
$ws.Range("S1221").Value = 22
$ws.Range("T1221").Value = 2
$ws.Range("U1221").Value = 11
$ws.Range("V1221").Value = 4

# Row 1222
$ws.Range("A1222").Value = "Entrainement"
$ws.Range("B1222").Value = 46050
$ws.Range("C1222").Value = "Global"
$ws.Range("D1222").Value = "J-3"
$ws.Range("E1222").Value = "Mattheo Haon"
$ws.Range("F1222").Value = "right back"
$ws.Range("G1222").Value = "01:41:05"
$ws.Range("H1222").Value = 7.38
$ws.Range("I1222").Value = 0.83
$ws.Range("J1222").Value = 6.55
$ws.Range("K1222").Value = 0.55
$ws.Range("L1222").Value = 0.18
$ws.Range("M1222").Value = 0.1
$ws.Range("N1222").Value = 0.01
$ws.Range("O1222").Value = 7
$ws.Range("P1222").Value = 4.32
$ws.Range("Q1222").Value = 30.88
$ws.Range("R1222").Value = 4.22
$ws.Range("S1222").Value = 20
$ws.Range("T1222").Value = 5
$ws.Range("U1222").Value = 16
$ws.Range("V1222").Value = 3

# Row 1223
$ws.Range("A1223").Value = "Entrainement"
$ws.Range("B1223").Value = 46050
$ws.Range("C1223").Value = "Global"
$ws.Range("D1223").Value = "J-3"
$ws.Range("E1223").Value = "Mehdi Boussaid"
$ws.Range("F1223").Value = "center midfield"
$ws.Range("G1223").Value = "01:40:30"
$ws.Range("H1223").Value = 7.67
$ws.Range("I1223").Value = 1.08
$ws.Range("J1223").Value = 6.57
$ws.Range("K1223").Value = 0.8
$ws.Range("L1223").Value = 0.25
$ws.Range("M1223").Value = 0.05
$ws.Range("N1223").Value = 0
$ws.Range("O1223").Value = 5
$ws.Range("P1223").Value = 4.43
$ws.Range("Q1223").Value = 28.5
$ws.Range("R1223").Value = 4.24
$ws.Range("S1223").Value = 36
$ws.Range("T1223").Value = 2
$ws.Range("U1223").Value = 20
$ws.Range("V1223").Value = 4

# Row 1224
$ws.Range("A1224").Value = "Entrainement"
$ws.Range("B1224").Value = 46050
$ws.Range("C1224").Value = "Global"
$ws.Range("D1224").Value = "J-3"
$ws.Range("E1224").Value = "Kamal Bafounta"
$ws.Range("F1224").Value = "center midfield"
$ws.Range("G1224").Value = "01:37:24"
$ws.Range("H1224").Value = 8.6
$ws.Range("I1224").Value = 1.49
$ws.Range("J1224").Value = 7.09
$ws.Range("K1224").Value = 0.96
$ws.Range("L1224").Value = 0.37
$ws.Range("M1224").Value = 0.18
$ws.Range("N1224").Value = 0
$ws.Range("O1224").Value = 13
$ws.Range("P1224").Value = 5.21
$ws.Range("Q1224").Value = 29.78
$ws.Range("R1224").Value = 4.88
$ws.Range("S1224").Value = 47
$ws.Range("T1224").Value = 4
$ws.Range("U1224").Value = 36
$ws.Range("V1224").Value = 3

# Row 1225
$ws.Range("A1225").Value = "Entrainement"
$ws.Range("B1225").Value = 46050
$ws.Range("C1225").Value = "Global"
$ws.Range("D1225").Value = "J-3"
$ws.Range("E1225").Value = "Hedi Nasri"
$ws.Range("F1225").Value = "right back"
$ws.Range("G1225").Value = "01:38:43"
$ws.Range("H1225").Value = 7.14
$ws.Range("I1225").Value = 0.94
$ws.Range("J1225").Value = 6.18
$ws.Range("K1225").Value = 0.48
$ws.Range("L1225").Value = 0.22
$ws.Range("M1225").Value = 0.2
$ws.Range("N1225").Value = 0.07
$ws.Range("O1225").Value = 16
$ws.Range("P1225").Value = 4.09
$ws.Range("Q1225").Value = 33.13
$ws.Range("R1225").Value = 5.74
$ws.Range("S1225").Value = 40
$ws.Range("T1225").Value = 25
$ws.Range("U1225").Value = 42
$ws.Range("V1225").Value = 8

# Row 1226
$ws.Range("A1226").Value = "Entrainement"
$ws.Range("B1226").Value = 46050
$ws.Range("C1226").Value = "Global"
$ws.Range("D1226").Value = "J-3"
$ws.Range("E1226").Value = "Amine Taiar"
$ws.Range("F1226").Value = "center back"
$ws.Range("G1226").Value = "01:40:22"
$ws.Range("H1226").Value = 7.45
$ws.Range("I1226").Value = 0.66
$ws.Range("J1226").Value = 6.75
$ws.Range("K1226").Value = 0.49
$ws.Range("L1226").Value = 0.16
$ws.Range("M1226").Value = 0.05
$ws.Range("N1226").Value = 0
$ws.Range("O1226").Value = 11
$ws.Range("P1226").Value = 4.05
$ws.Range("Q1226").Value = 29.32
$ws.Range("R1226").Value = 6.51
$ws.Range("S1226").Value = 226
$ws.Range("T1226").Value = 141
$ws.Range("U1226").Value = 175
$ws.Range("V1226").Value = 119

# Row 1227
$ws.Range("A1227").Value = "Entrainement"
$ws.Range("B1227").Value = 46050
$ws.Range("C1227").Value = "Global"
$ws.Range("D1227").Value = "J-3"
$ws.Range("E1227").Value = "Theo Owono"
$ws.Range("F1227").Value = "center midfield"
$ws.Range("G1227").Value = "01:39:10"
$ws.Range("H1227").Value = 7.09
$ws.Range("I1227").Value = 0.88
$ws.Range("J1227").Value = 6.2
$ws.Range("K1227").Value = 0.65
$ws.Range("L1227").Value = 0.23
$ws.Range("M1227").Value = 0.02
$ws.Range("N1227").Value = 0
$ws.Range("O1227").Value = 2
$ws.Range("P1227").Value = 4.27
$ws.Range("Q1227").Value = 26.67
$ws.Range("R1227").Value = 4.13
$ws.Range("S1227").Value = 33
$ws.Range("T1227").Value = 1
$ws.Range("U1227").Value = 27
$ws.Range("V1227").Value = 7

# Row 1228
$ws.Range("A1228").Value = "Entrainement"
$ws.Range("B1228").Value = 46050
$ws.Range("C1228").Value = "Global"
$ws.Range("D1228").Value = "J-3"
$ws.Range("E1228").Value = "Malik Boussaid"
$ws.Range("F1228").Value = "right back"
$ws.Range("G1228").Value = "01:41:05"
$ws.Range("H1228").Value = 7.75
$ws.Range("I1228").Value = 0.95
$ws.Range("J1228").Value = 6.78
$ws.Range("K1228").Value = 0.53
$ws.Range("L1228").Value = 0.28
$ws.Range("M1228").Value = 0.12
$ws.Range("N1228").Value = 0.03
$ws.Range("O1228").Value = 12
$ws.Range("P1228").Value = 4.28
$ws.Range("Q1228").Value = 32.81
$ws.Range("R1228").Value = 4.6
$ws.Range("S1228").Value = 62
$ws.Range("T1228").Value = 6
$ws.Range("U1228").Value = 43
$ws.Range("V1228").Value = 14

# Row 1229
$ws.Range("A1229").Value = "Entrainement"
$ws.Range("B1229").Value = 46050
$ws.Range("C1229").Value = "Global"
$ws.Range("D1229").Value = "J-3"
$ws.Range("E1229").Value = "Jeremie Laurent"
$ws.Range("F1229").Value = "left forward"
$ws.Range("G1229").Value = "00:49:35"
$ws.Range("H1229").Value = 2.88
$ws.Range("I1229").Value = 0.37
$ws.Range("J1229").Value = 2.51
$ws.Range("K1229").Value = 0.11
$ws.Range("L1229").Value = 0.12
$ws.Range("M1229").Value = 0.13
$ws.Range("N1229").Value = 0.01
$ws.Range("O1229").Value = 8
$ws.Range("P1229").Value = 3.4
$ws.Range("Q1229").Value = 30.81
$ws.Range("R1229").Value = 4.99
$ws.Range("S1229").Value = 13
$ws.Range("T1229").Value = 7
$ws.Range("U1229").Value = 9
$ws.Range("V1229").Value = 1

# Update the sheet view to match where the author ended up after entering the new rows:
# scrolled further down and landed the selection on E1232 (one row below the new data + a blank row).
[void]$ws.Range("A1197").Select()
[void]$ws.Range("E1232").Select()
